$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell B3 ("SamplesTab" query): narrow the sample_type filter from the
# IN-list ('Metastatic','Blood Derived Normal') down to just 'Metastatic',
# matching the commit's regression-suite intent. Do it as a targeted
# substring replace against the live cell value so the rest of the query
# text/formatting is left untouched.
$oldQueryCell = $ws.Range("B3")
$currentQuery = $oldQueryCell.Value()

$oldClause = "   s.phs_accession = 'phs003155' AND smp.sample_type IN ('Metastatic','Blood Derived Normal')"
$newClause = "  s.phs_accession = 'phs003155' AND smp.sample_type = 'Metastatic'"

$updatedQuery = $currentQuery.Replace($oldClause, $newClause)
$oldQueryCell.Value = $updatedQuery

# --- View state: the worksheet was scrolled/reselected before saving
# (topLeftCell A3 -> A4, selection B3 -> C15). The row scrolled by one;
# the left-most visible column (A) is unchanged.
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("C15").Select()
